$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3148
$ws1.Range("F3").Value = 540
$ws1.Range("F4").Value = 1104
$ws1.Range("F5").Value = 92
$ws1.Range("F6").Value = 46
$ws1.Range("F9").Value = 1134
$ws1.Range("F10").Value = 15887
$ws1.Range("F11").Value = 253
$ws1.Range("F14").Value = 6228
$ws1.Range("F15").Value = 627
$ws1.Range("F16").Value = 115
$ws1.Range("F17").Value = 70
$ws1.Range("F18").Value = 11
$ws1.Range("F23").Value = 22
$ws1.Range("F24").Value = 16
$ws1.Range("F26").Value = 214
$ws1.Range("F27").Value = 877
$ws1.Range("F29").Value = 5011
$ws1.Range("F30").Value = 492
$ws1.Range("F31").Value = 11132
$ws1.Range("F32").Value = 1237
$ws1.Range("F34").Value = 132
$ws1.Range("F35").Value = 185

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3148
$ws4.Range("F4").Value = 540
$ws4.Range("F5").Value = 1104
$ws4.Range("F6").Value = 92
$ws4.Range("F7").Value = 46
$ws4.Range("F10").Value = 1134
$ws4.Range("F11").Value = 15887
$ws4.Range("F12").Value = 253
$ws4.Range("F15").Value = 6228
$ws4.Range("F16").Value = 627
$ws4.Range("F17").Value = 115
$ws4.Range("F18").Value = 70
$ws4.Range("F19").Value = 11
$ws4.Range("F24").Value = 22
$ws4.Range("F25").Value = 16
$ws4.Range("F27").Value = 214
$ws4.Range("F28").Value = 877
$ws4.Range("F30").Value = 5011
$ws4.Range("F31").Value = 492
$ws4.Range("F33").Value = 11132
$ws4.Range("F34").Value = 1237
$ws4.Range("F36").Value = 132
$ws4.Range("F37").Value = 185
